# Insert 4 new price rows above row 948 (weekly Femacal de La Calera - Limon update),
# shifting the existing rows 948:998 down to 952:1002.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("948:951").Insert()

# Shared values for the 4 new rows (same market / region / product context as the
# surrounding rows).
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"
$variedad  = "Sin especificar"
$unidad    = "$/malla 16 kilos"
$origen    = "Provincia de Quillota"
$kgUnidad  = 16
$fecha     = 44706

# Row 948: 1a amarillo
$ws.Cells.Item(948, 1).Value  = $mercadoId
$ws.Cells.Item(948, 2).Value  = $mercado
$ws.Cells.Item(948, 3).Value  = $region
$ws.Cells.Item(948, 4).Value  = $fecha
$ws.Cells.Item(948, 5).Value  = $codreg
$ws.Cells.Item(948, 6).Value  = $tipo
$ws.Cells.Item(948, 7).Value  = $productoId
$ws.Cells.Item(948, 8).Value  = $producto
$ws.Cells.Item(948, 9).Value  = $categoriaId
$ws.Cells.Item(948, 10).Value = $categoria
$ws.Cells.Item(948, 11).Value = $variedad
$ws.Cells.Item(948, 12).Value = "1a amarillo"
$ws.Cells.Item(948, 13).Value = 100
$ws.Cells.Item(948, 14).Value = 6000
$ws.Cells.Item(948, 15).Value = 6000
$ws.Cells.Item(948, 16).Value = 6000
$ws.Cells.Item(948, 17).Value = $unidad
$ws.Cells.Item(948, 18).Value = $origen
$ws.Cells.Item(948, 19).Value = 375
$ws.Cells.Item(948, 20).Value = $kgUnidad

# Row 949: 1a plateado
$ws.Cells.Item(949, 1).Value  = $mercadoId
$ws.Cells.Item(949, 2).Value  = $mercado
$ws.Cells.Item(949, 3).Value  = $region
$ws.Cells.Item(949, 4).Value  = $fecha
$ws.Cells.Item(949, 5).Value  = $codreg
$ws.Cells.Item(949, 6).Value  = $tipo
$ws.Cells.Item(949, 7).Value  = $productoId
$ws.Cells.Item(949, 8).Value  = $producto
$ws.Cells.Item(949, 9).Value  = $categoriaId
$ws.Cells.Item(949, 10).Value = $categoria
$ws.Cells.Item(949, 11).Value = $variedad
$ws.Cells.Item(949, 12).Value = "1a plateado"
$ws.Cells.Item(949, 13).Value = 198
$ws.Cells.Item(949, 14).Value = 6000
$ws.Cells.Item(949, 15).Value = 6500
$ws.Cells.Item(949, 16).Value = 6247
$ws.Cells.Item(949, 17).Value = $unidad
$ws.Cells.Item(949, 18).Value = $origen
$ws.Cells.Item(949, 19).Value = 390
$ws.Cells.Item(949, 20).Value = $kgUnidad

# Row 950: 2a amarillo
$ws.Cells.Item(950, 1).Value  = $mercadoId
$ws.Cells.Item(950, 2).Value  = $mercado
$ws.Cells.Item(950, 3).Value  = $region
$ws.Cells.Item(950, 4).Value  = $fecha
$ws.Cells.Item(950, 5).Value  = $codreg
$ws.Cells.Item(950, 6).Value  = $tipo
$ws.Cells.Item(950, 7).Value  = $productoId
$ws.Cells.Item(950, 8).Value  = $producto
$ws.Cells.Item(950, 9).Value  = $categoriaId
$ws.Cells.Item(950, 10).Value = $categoria
$ws.Cells.Item(950, 11).Value = $variedad
$ws.Cells.Item(950, 12).Value = "2a amarillo"
$ws.Cells.Item(950, 13).Value = 100
$ws.Cells.Item(950, 14).Value = 4500
$ws.Cells.Item(950, 15).Value = 4500
$ws.Cells.Item(950, 16).Value = 4500
$ws.Cells.Item(950, 17).Value = $unidad
$ws.Cells.Item(950, 18).Value = $origen
$ws.Cells.Item(950, 19).Value = 281
$ws.Cells.Item(950, 20).Value = $kgUnidad

# Row 951: 2a plateado
$ws.Cells.Item(951, 1).Value  = $mercadoId
$ws.Cells.Item(951, 2).Value  = $mercado
$ws.Cells.Item(951, 3).Value  = $region
$ws.Cells.Item(951, 4).Value  = $fecha
$ws.Cells.Item(951, 5).Value  = $codreg
$ws.Cells.Item(951, 6).Value  = $tipo
$ws.Cells.Item(951, 7).Value  = $productoId
$ws.Cells.Item(951, 8).Value  = $producto
$ws.Cells.Item(951, 9).Value  = $categoriaId
$ws.Cells.Item(951, 10).Value = $categoria
$ws.Cells.Item(951, 11).Value = $variedad
$ws.Cells.Item(951, 12).Value = "2a plateado"
$ws.Cells.Item(951, 13).Value = 193
$ws.Cells.Item(951, 14).Value = 4500
$ws.Cells.Item(951, 15).Value = 5000
$ws.Cells.Item(951, 16).Value = 4754
$ws.Cells.Item(951, 17).Value = $unidad
$ws.Cells.Item(951, 18).Value = $origen
$ws.Cells.Item(951, 19).Value = 297
$ws.Cells.Item(951, 20).Value = $kgUnidad
